$d = $word.ActiveDocument

$replacements = @(
    @("15×61=915", "15×64=960"),
    @("53×91=4823", "84×34=2856"),
    @("78×44=3432", "19×76=1444"),
    @("69×58=4002", "60×91=5460"),
    @("11×67=737", "47×21=987"),
    @("32×21=672", "70×51=3570"),
    @("25×68=1700", "21×76=1596"),
    @("14×57=798", "66×80=5280"),
    @("11×42=462", "56×91=5096"),
    @("37×53=1961", "25×39=975"),
    @("51×92=4692", "81×56=4536"),
    @("37×84=3108", "52×78=4056"),
    @("96×92=8832", "87×47=4089"),
    @("15×68=1020", "36×63=2268"),
    @("33×50=1650", "50×59=2950"),
    @("88×52=4576", "40×21=840"),
    @("85×96=8160", "22×52=1144"),
    @("64×13=832", "23×90=2070"),
    @("33×17=561", "91×16=1456"),
    @("46×49=2254", "16×56=896"),
    @("80×65=5200", "34×89=3026"),
    @("59×25=1475", "60×86=5160"),
    @("68×60=4080", "78×77=6006"),
    @("54×47=2538", "71×54=3834"),
    @("37×74=2738", "64×36=2304")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
